$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AS1").Value = "Week 44"
$ws.Range("AS2").Value = 2.67
$ws.Range("AS4").Value = 2
$ws.Range("AS5").Value = 5
$ws.Range("AS6").Value = 5.5
$ws.Range("AS10").Value = 10

[void]$ws.Range("AQ12").Select()
